$wb = $excel.ActiveWorkbook

# GitHub blob URL prefix used by the existing hyperlinks in this workbook.
$baseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/d259aa2299996133a837809225144a69a6ee9e85"

$newFile = "70485d62-d35f-4af0-aeeb-4b07495edfce.md"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview" -- File Name / zh-cn / de-de summary (A1:C4 -> A1:C5)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert a new row 3 - this pushes the old rows 3 & 4 down to 4 & 5 and
# carries the hyperlink-style formatting (s="1") on column A along with it.
$ws1.Rows(3).Insert()

$ws1.Range("A3").Value = $newFile
$ws1.Range("B3").Value = "Handoff transform failed"
$ws1.Range("C3").Value = "Handoff transform failed"

# Rebuild the hyperlinks for column A in row order so the relationship ids
# line up sequentially (rId2..rId5) the way Excel lays them out.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $baseUrl + "/e2e/08e16aaf-fb5d-4c93-a634-363ff4d4e53a.md", "", "", "08e16aaf-fb5d-4c93-a634-363ff4d4e53a.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), $baseUrl + "/e2e/" + $newFile, "", "", $newFile)
$ws1.Hyperlinks.Add($ws1.Range("A4"), $baseUrl + "/e2e/9528511d-8529-4a11-b078-0148a84d2a87.md", "", "", "9528511d-8529-4a11-b078-0148a84d2a87.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), $baseUrl + "/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn" detail sheet (A1:I4 -> A1:I5)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Rows(3).Insert()

$ws2.Range("A3").Value = $newFile
$ws2.Range("B3").Value = "Handoff transform failed"
$ws2.Range("D3").Value = "0001-01-01 00:00:00"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Ignored"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $baseUrl + "/e2e/08e16aaf-fb5d-4c93-a634-363ff4d4e53a.md", "", "", "08e16aaf-fb5d-4c93-a634-363ff4d4e53a.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), $baseUrl + "/e2e/" + $newFile, "", "", $newFile)
$ws2.Hyperlinks.Add($ws2.Range("A4"), $baseUrl + "/e2e/9528511d-8529-4a11-b078-0148a84d2a87.md", "", "", "9528511d-8529-4a11-b078-0148a84d2a87.md")
$ws2.Hyperlinks.Add($ws2.Range("A5"), $baseUrl + "/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet 3: "de-de" detail sheet (A1:I4 -> A1:I5) -- mirrors sheet 2
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Rows(3).Insert()

$ws3.Range("A3").Value = $newFile
$ws3.Range("B3").Value = "Handoff transform failed"
$ws3.Range("D3").Value = "0001-01-01 00:00:00"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Ignored"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $baseUrl + "/e2e/08e16aaf-fb5d-4c93-a634-363ff4d4e53a.md", "", "", "08e16aaf-fb5d-4c93-a634-363ff4d4e53a.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), $baseUrl + "/e2e/" + $newFile, "", "", $newFile)
$ws3.Hyperlinks.Add($ws3.Range("A4"), $baseUrl + "/e2e/9528511d-8529-4a11-b078-0148a84d2a87.md", "", "", "9528511d-8529-4a11-b078-0148a84d2a87.md")
$ws3.Hyperlinks.Add($ws3.Range("A5"), $baseUrl + "/.localization-config", "", "", ".localization-config")
